$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B54 was stored as a text "4"; correct it to a real number 4.
$ws.Range("B54").Value = 4

# Append a new annotation row (row 55).
$ws.Range("A55").Value = "Ruilin"

# B55's "3" stays a text value (matches source data pattern), not a number.
# Force text entry, then drop back to the default "Normal" style so no
# stray number-format styling is left attached to the cell.
$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = "3"
$ws.Range("B55").Style = "Normal"

$ws.Range("C55").Value = "无"
$ws.Range("D55").Value = "DIS"
$ws.Range("E55").Value = "MET"
$ws.Range("F55").Value = "0ffe4b07-d72b-4753-8576-ca80ee89bdb3"
$ws.Range("G55").Value = "SJzMATlAZ_annotated.xlsx"
$ws.Range("H55").Value = "We avoid using k-means because it requires knowing the number of clusters a priory."
